$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.268.09"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.689.92"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'219.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "'0.5246"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.13%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").Value = "'0.2695"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.86%  "
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D10").Value = "'22.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.35%  "
$ws.Range("D11").Value = "'0.07461"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("D12").Value = "1.687.54"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Value = "'0.5858"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").Value = "'0.000008548"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").Value = "'64.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("D17").Value = "26.285.19"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "'4.968"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").Value = "'190.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("D22").Value = "'6.239"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "'145.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("D25").Value = "'7.680"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").Value = "'0.1241"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.99%  "
$ws.Range("D27").Value = "'15.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").Value = "'0.06671"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +13.39%  "
$ws.Range("D29").Value = "'1.353"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.73%  "
$ws.Range("D30").Value = "'1.331"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("D31").Value = "'3.598"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.63%  "
$ws.Range("E32").Value = "  +1.34%  "
$ws.Range("D33").Value = "'1.668"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("D35").Value = "'0.6205"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.19%  "
$ws.Range("D36").Value = "'2.388"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("D37").Value = "'2.713"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.57%  "
$ws.Range("D38").Value = "'6.281"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.13%  "
$ws.Range("E39").Value = "  +0.73%  "
$ws.Range("D40").Value = "1.103.40"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("D41").Value = "'0.8832"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.81%  "
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("D43").Value = "'100.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("D44").Value = "1.837.91"
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("D45").Value = "'0.00000000116"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.53%  "
$ws.Range("D46").Value = "'56.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.63%  "
$ws.Range("D47").Value = "'1.008"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("D48").Value = "'8.150"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D49").Value = "'0.05263"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.57%  "
$ws.Range("D50").Value = "'0.4296"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").Value = "'6.031"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.79%  "
